# SCD0018-019: update TC_ID from "DGS-311" to "SCD0018-019" and rename the
# sheet from "SCD0296" to "SCD0018" to match the new test-case id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the new TC id prefix.
$ws.Name = "SCD0018"

# The TC_ID column (B) held "DGS-311" for every data row (2-6); replace it
# with the new case id "SCD0018-019".
$ws.Range("B2:B6").Value = "SCD0018-019"

# The TC_ID column is a bit wider now that it holds the longer id - resize
# it to fit the new content, like the author would have after editing.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Leave the cursor on the scenario description cell that was being edited.
$ws.Range("C3").Select()
